# The workbook originally has a single sheet, "String".
# This adds a second sheet, "Dummy", positioned right after "String",
# containing a single cell (B1) with the text "col3 of 2".

$wb = $excel.ActiveWorkbook

$stringSheet = $wb.Worksheets.Item(1)

# Insert the new sheet immediately after the "String" sheet.
$dummySheet = $wb.Worksheets.Add($null, $stringSheet)
$dummySheet.Name = "Dummy"
$dummySheet.Range("B1").Value = "col3 of 2"

# Keep "String" as the active/selected sheet, matching the original workbook.
$stringSheet.Activate()
